$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 920.8570999999999
$ws.Range("I103").Value = 850
$ws.Range("J103").Value = 932.6667
$ws.Range("K103").Value = 2550
$ws.Range("L103").Value = 2798.0001
$ws.Range("M103").Value = -1964
$ws.Range("N103").Value = -3970.0001

$ws.Range("H135").Value = 1269.25
$ws.Range("I135").Value = 1269.25
$ws.Range("K135").Value = 11423.25
$ws.Range("M135").Value = -8888.25

$ws.Range("H137").Value = 5039.467
$ws.Range("I137").Value = 2849.5
$ws.Range("J137").Value = 5835.8184
$ws.Range("K137").Value = 8548.5
$ws.Range("L137").Value = 17507.4552
$ws.Range("M137").Value = -5998.5
$ws.Range("N137").Value = -22607.4552

$ws.Range("H138").Value = 2637.8
$ws.Range("I138").Value = 1068.5161
$ws.Range("J138").Value = 4664.7915
$ws.Range("K138").Value = 3205.5483
$ws.Range("L138").Value = 13994.3745
$ws.Range("M138").Value = 1934.4517
$ws.Range("N138").Value = -24274.3745

$ws.Range("H141").Value = 2587.5715
$ws.Range("I141").Value = 1871.6666
$ws.Range("K141").Value = 5614.9998
$ws.Range("M141").Value = -434.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4411.364
$ws.Range("I32").Value = 3728.5806
$ws.Range("J32").Value = 14994.5
$ws.Range("K32").Value = 3728.5806
$ws.Range("L32").Value = 14994.5
$ws.Range("M32").Value = -3441.5806
$ws.Range("N32").Value = -15568.5

$ws.Range("H61").Value = 2827.2307
$ws.Range("I61").Value = 2796.182
$ws.Range("J61").Value = 2998
$ws.Range("K61").Value = 2796.182
$ws.Range("L61").Value = 2998
$ws.Range("M61").Value = -2584.182
$ws.Range("N61").Value = -3422

$ws.Range("H74").Value = 1907.6
$ws.Range("I74").Value = 2084.75
$ws.Range("K74").Value = 2084.75
$ws.Range("M74").Value = -1210.75

$ws.Range("H77").Value = 1907.6
$ws.Range("I77").Value = 2084.75
$ws.Range("K77").Value = 10423.75
$ws.Range("M77").Value = -6055.75

$ws.Range("H132").Value = 2123.9736
$ws.Range("I132").Value = 1960.3429
$ws.Range("K132").Value = 5881.028700000001
$ws.Range("M132").Value = -3351.028700000001

$ws.Range("H136").Value = 2827.2307
$ws.Range("I136").Value = 2796.182
$ws.Range("J136").Value = 2998
$ws.Range("K136").Value = 8388.545999999998
$ws.Range("L136").Value = 8994
$ws.Range("M136").Value = -5838.545999999998
$ws.Range("N136").Value = -14094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3623.0667
$ws.Range("I134").Value = 3703.3572
$ws.Range("K134").Value = 11110.0716
$ws.Range("M134").Value = -8575.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2994.75
$ws.Range("I16").Value = 2994.75
$ws.Range("K16").Value = 2994.75
$ws.Range("M16").Value = -2707.75

$ws.Range("H31").Value = 2130.889
$ws.Range("I31").Value = 2285.125
$ws.Range("K31").Value = 2285.125
$ws.Range("M31").Value = -1990.125

$ws.Range("H34").Value = 2130.889
$ws.Range("I34").Value = 2285.125
$ws.Range("K34").Value = 2285.125
$ws.Range("M34").Value = -2083.125

$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = $null
$ws.Range("N69").Value = -31498

$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = $null
$ws.Range("N72").Value = -97488

$ws.Range("H86").Value = 18436.264
$ws.Range("J86").Value = 37040.375
$ws.Range("L86").Value = 37040.375
$ws.Range("N86").Value = -39286.375

$ws.Range("H89").Value = 18436.264
$ws.Range("J89").Value = 37040.375
$ws.Range("L89").Value = 185201.875
$ws.Range("N89").Value = -196433.875

$ws.Range("H99").Value = 4329.6665
$ws.Range("I99").Value = 4499.5
$ws.Range("K99").Value = 4499.5
$ws.Range("M99").Value = -3001.5

$ws.Range("H113").Value = 2994.75
$ws.Range("I113").Value = 2994.75
$ws.Range("K113").Value = 2994.75
$ws.Range("M113").Value = -824.75

$ws.Range("H126").Value = 4329.6665
$ws.Range("I126").Value = 4499.5
$ws.Range("K126").Value = 13498.5
$ws.Range("M126").Value = -11028.5

$ws.Range("H134").Value = 5650.5557
$ws.Range("I134").Value = 5650.5557
$ws.Range("K134").Value = 16951.6671
$ws.Range("M134").Value = -14416.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 1533.3334
$ws.Range("I104").Value = 1200
$ws.Range("K104").Value = 3600
$ws.Range("M104").Value = -979

$ws.Range("H132").Value = 3547.8
$ws.Range("I132").Value = 3547.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 31930.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -29400.2
$ws.Range("N132").Value = $null

$ws.Range("H139").Value = 3332
$ws.Range("I139").Value = 1930
$ws.Range("J139").Value = 4033
$ws.Range("K139").Value = 5790
$ws.Range("L139").Value = 12099
$ws.Range("M139").Value = -650
$ws.Range("N139").Value = -22379

$ws.Range("H140").Value = 3241.5
$ws.Range("I140").Value = 3241.5
$ws.Range("K140").Value = 9724.5
$ws.Range("M140").Value = -4544.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 799
$ws.Range("I122").Value = 799
$ws.Range("K122").Value = 2397
$ws.Range("M122").Value = 53

$ws.Range("H132").Value = 2999.8333
$ws.Range("I132").Value = 2999.8
$ws.Range("K132").Value = 8999.400000000001
$ws.Range("M132").Value = -6469.400000000001

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2387.375
$ws.Range("I7").Value = 2442.7144
$ws.Range("K7").Value = 2442.7144
$ws.Range("M7").Value = -2330.7144

$ws.Range("H40").Value = 4998
$ws.Range("I40").Value = 4797.4
$ws.Range("K40").Value = 4797.4
$ws.Range("M40").Value = -4661.4

$ws.Range("H122").Value = 3199.9285
$ws.Range("I122").Value = 2710.182
$ws.Range("K122").Value = 8130.545999999999
$ws.Range("M122").Value = -5680.545999999999

$ws.Range("H126").Value = 2387.375
$ws.Range("I126").Value = 2442.7144
$ws.Range("K126").Value = 7328.1432
$ws.Range("M126").Value = -4858.1432

$ws.Range("H132").Value = 2470.9443
$ws.Range("I132").Value = 1900.1
$ws.Range("J132").Value = 3184.5
$ws.Range("K132").Value = 5700.299999999999
$ws.Range("L132").Value = 9553.5
$ws.Range("M132").Value = -3170.299999999999
$ws.Range("N132").Value = -14613.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10277.777
$ws.Range("J4").Value = 10312.25
$ws.Range("L4").Value = 10312.25
$ws.Range("N4").Value = -10538.25

$ws.Range("H5").Value = 9000
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9224

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = $null

$ws.Range("H136").Value = 1058.6072
$ws.Range("I136").Value = 883.9583
$ws.Range("K136").Value = 2651.8749
$ws.Range("M136").Value = -101.8748999999998

Write-Output "Updated 39 rows across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
